$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Total" (column B) and "Community" (column D) monthly consumption
# values fixing the big mistake in the original data.

$ws.Range("B2").Value = 7873.70901350002
$ws.Range("D2").Value = 505.2825002833334

$ws.Range("B3").Value = 7366.369618183351
$ws.Range("D3").Value = 479.8452055333333

$ws.Range("B4").Value = 7900.730378083354
$ws.Range("D4").Value = 529.6545940666667

$ws.Range("B5").Value = 7617.880432633352
$ws.Range("D5").Value = 495.8968661333333

$ws.Range("B6").Value = 7903.597525466686
$ws.Range("D6").Value = 523.54358935

$ws.Range("B7").Value = 7655.652308883353
$ws.Range("D7").Value = 509.1682569166667

$ws.Range("B8").Value = 7895.98543095002
$ws.Range("D8").Value = 518.8358951499999

$ws.Range("B9").Value = 7895.23621675002
$ws.Range("D9").Value = 523.1529089666667

$ws.Range("B10").Value = 7664.542564450019
$ws.Range("D10").Value = 500.90419505

$ws.Range("B11").Value = 7893.204864216686
$ws.Range("D11").Value = 523.8685568833333

$ws.Range("B12").Value = 7666.541936400019
$ws.Range("D12").Value = 492.1805045666667

$ws.Range("B13").Value = 7641.802333766685
$ws.Range("D13").Value = 504.9592266333333
